$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before D (shifts existing D:K to F:M)
$ws.Columns("D:E").Insert()

# Copy number formatting from column F into new D:E columns for each data block
$ws.Range("F7:F35").Copy()
$ws.Range("D7:E35").PasteSpecial(-4122)
$ws.Range("F38:F77").Copy()
$ws.Range("D38:E77").PasteSpecial(-4122)
$ws.Range("F80:F102").Copy()
$ws.Range("D80:E102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate new column D and E values
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("D8").Value = 1260300
$ws.Range("E8").Value = 1186200
$ws.Range("D9").Value = 861200
$ws.Range("E9").Value = 820700
$ws.Range("D10").Value = 399100
$ws.Range("E10").Value = 365500
$ws.Range("D12").Value = "NA"
$ws.Range("E12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = 25500
$ws.Range("E14").Value = 6600
$ws.Range("D15").Value = 4800
$ws.Range("E15").Value = 3600
$ws.Range("D17").Value = 1094900
$ws.Range("E17").Value = 1023000
$ws.Range("D18").Value = 165400
$ws.Range("E18").Value = 163200
$ws.Range("D20").Value = 2300
$ws.Range("E20").Value = -6300
$ws.Range("D21").Value = 200500
$ws.Range("E21").Value = 189600
$ws.Range("D22").Value = 49200
$ws.Range("E22").Value = 47900
$ws.Range("D23").Value = 118500
$ws.Range("E23").Value = 109000
$ws.Range("D24").Value = -14500
$ws.Range("E24").Value = 33400
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = 133000
$ws.Range("E26").Value = 75600
$ws.Range("D27").Value = 132500
$ws.Range("E27").Value = 75200
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = 67300
$ws.Range("E29").Value = 3400
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = -2300
$ws.Range("E32").Value = 6300
$ws.Range("D33").Value = 199800
$ws.Range("E33").Value = 78600
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = 199800
$ws.Range("E35").Value = 78600
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("D41").Value = 271700
$ws.Range("E41").Value = 191300
$ws.Range("D42").Value = 0
$ws.Range("E42").Value = 0
$ws.Range("D43").Value = 613100
$ws.Range("E43").Value = 630200
$ws.Range("D44").Value = 544900
$ws.Range("E44").Value = 605400
$ws.Range("D45").Value = 125100
$ws.Range("E45").Value = 168500
$ws.Range("D46").Value = 1554800
$ws.Range("E46").Value = 1595400
$ws.Range("D47").Value = 0
$ws.Range("E47").Value = 0
$ws.Range("D48").Value = 1036200
$ws.Range("E48").Value = 1022000
$ws.Range("D49").Value = 2049300
$ws.Range("E49").Value = 2054200
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 409900
$ws.Range("E52").Value = 325400
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 5050200
$ws.Range("E54").Value = 4997000
$ws.Range("D57").Value = 765000
$ws.Range("E57").Value = 775300
$ws.Range("D58").Value = 237700
$ws.Range("E58").Value = 313400
$ws.Range("D59").Value = 485900
$ws.Range("E59").Value = 477900
$ws.Range("D60").Value = 1488600
$ws.Range("E60").Value = 1566600
$ws.Range("D61").Value = 3236500
$ws.Range("E61").Value = 3242500
$ws.Range("D62").Value = 673700
$ws.Range("E62").Value = 633600
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 5398800
$ws.Range("E66").Value = 5442700
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = 1835500
$ws.Range("E72").Value = 1646700
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = -348600
$ws.Range("E76").Value = -445700
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("D81").Value = 199800
$ws.Range("E81").Value = 78600
$ws.Range("D83").Value = 32800
$ws.Range("E83").Value = 32700
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = 278000
$ws.Range("E89").Value = 113400
$ws.Range("D91").Value = -53800
$ws.Range("E91").Value = -41100
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = -55500
$ws.Range("E94").Value = -130400
$ws.Range("D96").Value = -24800
$ws.Range("E96").Value = -25300
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = -144100
$ws.Range("E100").Value = 28500
$ws.Range("D101").Value = 2000
$ws.Range("E101").Value = -300
$ws.Range("D102").Value = 80400
$ws.Range("E102").Value = 11200

# Fix cells whose shifted values were restated with different figures
$ws.Range("H9").Value = 858500
$ws.Range("H10").Value = 369300
$ws.Range("H17").Value = 1090700
$ws.Range("H18").Value = 137100
$ws.Range("F20").Value = 5200
$ws.Range("G20").Value = -6200
$ws.Range("H20").Value = 12400
$ws.Range("F21").Value = 198300
$ws.Range("G21").Value = 194100
$ws.Range("H21").Value = 184400
$ws.Range("F22").Value = 48600
$ws.Range("G22").Value = 47800
$ws.Range("H22").Value = 30500
$ws.Range("F32").Value = -5200
$ws.Range("G32").Value = 6200
$ws.Range("H32").Value = -12400
$ws.Range("I91").Value = -33300
$ws.Range("J91").Value = -42800
$ws.Range("H100").Value = -619100
$ws.Range("I100").Value = -1109100
